$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.325.67'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.010.44'
$ws.Range('E3').Value = '  +5.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.72'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5147'
$ws.Range('E7').Value = '  +1.93%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4262'
$ws.Range('E8').Value = '  +5.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08692'
$ws.Range('E9').Value = '  +5.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.22'
$ws.Range('E10').Value = '  +3.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.131'
$ws.Range('E11').Value = '  +2.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.72'
$ws.Range('E12').Value = '  +2.75%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.010.95'
$ws.Range('E13').Value = '  +5.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.566'
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.458'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.43'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001111'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06524'
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.85'
$ws.Range('E20').Value = '  +3.99%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.198'
$ws.Range('E22').Value = '  +4.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.393.21'
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.77'
$ws.Range('E24').Value = '  +4.14%  '
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.245.59'
$ws.Range('E26').Value = '  +5.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.34'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.62'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.408'
$ws.Range('E29').Value = '  +4.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '130.92'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.134'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1053'
$ws.Range('E32').Value = '  +1.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.052'
$ws.Range('E33').Value = '  +2.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.836'
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.362'
$ws.Range('E35').Value = '  +14.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02524'
$ws.Range('E36').Value = '  +3.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.456'
$ws.Range('E37').Value = '  +0.97%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06627'
$ws.Range('E38').Value = '  +3.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '12.31'
$ws.Range('E39').Value = '  +8.24%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2186'
$ws.Range('E40').Value = '  +1.97%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.002'
$ws.Range('E41').Value = '  +3.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6640'
$ws.Range('E42').Value = '  +2.68%  '
$ws.Range('E43').Value = '  +1.76%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.65'
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6148'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.179'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.662'
$ws.Range('E48').Value = '  +0.74%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.262'
$ws.Range('E49').Value = '  +4.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '123.99'
$ws.Range('E50').Value = '  +1.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.61'
$ws.Range('E51').Value = '  +2.43%  '
